# Update "想去人数" (F column) figures across the four sheets to reflect the
# latest scrape results (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1841
$ws1.Range("F4").Value  = 1501
$ws1.Range("F5").Value  = 862
$ws1.Range("F8").Value  = 13227
$ws1.Range("F9").Value  = 13093
$ws1.Range("F10").Value = 1002
$ws1.Range("F11").Value = 769
$ws1.Range("F16").Value = 2074
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 34
$ws1.Range("F19").Value = 49
$ws1.Range("F21").Value = 212

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 58
$ws2.Range("F9").Value = 16

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 191
$ws3.Range("F3").Value = 17

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 191
$ws4.Range("F3").Value  = 1841
$ws4.Range("F5").Value  = 1501
$ws4.Range("F6").Value  = 862
$ws4.Range("F10").Value = 13227
$ws4.Range("F11").Value = 13093
$ws4.Range("F12").Value = 1002
$ws4.Range("F13").Value = 769
$ws4.Range("F20").Value = 2074
$ws4.Range("F21").Value = 58
$ws4.Range("F22").Value = 34
$ws4.Range("F23").Value = 49
$ws4.Range("F26").Value = 58
$ws4.Range("F27").Value = 17
$ws4.Range("F28").Value = 212
$ws4.Range("F33").Value = 16
